# Insert a new paragraph right after the "Malloc to alloca" title
# paragraph (and before the "Didn't finish optimization..." paragraph),
# describing the results of the malloc-to-alloca optimization.

$d = $word.ActiveDocument

# The target paragraph ("Didn't finish optimization...") is currently
# paragraph #2; inserting a paragraph *before* it makes the new
# paragraph inherit that paragraph's (non-bold, non-centered) formatting,
# matching the style used throughout the rest of the document body.
$targetPara = $d.Paragraphs.Item(2)
$targetPara.Range.InsertParagraphBefore()

# The freshly inserted paragraph is now paragraph #2.
$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range

$fullText = "To replace malloc instruction to alloca, we need to go through all calls of malloc function and check whether it is freed ot no, after this we can replace instructions by alloca and remove free and malloc calls. The problem was appeared in replace instruction stage. Therefore, for now we need to find example of how to CreateAlloca of an array of pointers, also size check to prevent allocating huge memory on stack should be added. "

$newRange.Text = $fullText

$start = $newRange.Start

# Italicize the three inline code-style words/phrases that are emphasized
# in the source text: "alloca ", "free" and "malloc".
$run1 = $d.Range($start + 171, $start + 178)   # "alloca "
$run1.Italic = 1

$run2 = $d.Range($start + 189, $start + 193)   # "free"
$run2.Italic = 1

$run3 = $d.Range($start + 198, $start + 204)   # "malloc"
$run3.Italic = 1
